# Weekly update: Fruta / hortaliza, semanal
# Insert 3 new rows (18:20) for the new week's data, pushing the existing
# historical rows (old 18:44) down to 21:47. Then populate the 3 new rows
# with the latest Femacal de La Calera - Chirimoya price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old rows 18-44 down by three rows so the new week's three
# records can be inserted at the top of this date-ordered block.
$ws.Rows("18:20").Insert()

# --- Row 18: Especial ---
$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 44467
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100107
$ws.Cells.Item(18, 8).Value = "Otros"
$ws.Cells.Item(18, 9).Value = 100107002
$ws.Cells.Item(18, 10).Value = "Chirimoya"
$ws.Cells.Item(18, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(18, 12).Value = "Especial"
$ws.Cells.Item(18, 13).Value = 45
$ws.Cells.Item(18, 14).Value = 27000
$ws.Cells.Item(18, 15).Value = 27000
$ws.Cells.Item(18, 16).Value = 27000
$ws.Cells.Item(18, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(18, 19).Value = 2700
$ws.Cells.Item(18, 20).Value = 10

# --- Row 19: Primera ---
$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44467
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100107
$ws.Cells.Item(19, 8).Value = "Otros"
$ws.Cells.Item(19, 9).Value = 100107002
$ws.Cells.Item(19, 10).Value = "Chirimoya"
$ws.Cells.Item(19, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 40
$ws.Cells.Item(19, 14).Value = 25000
$ws.Cells.Item(19, 15).Value = 25000
$ws.Cells.Item(19, 16).Value = 25000
$ws.Cells.Item(19, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(19, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(19, 19).Value = 2500
$ws.Cells.Item(19, 20).Value = 10

# --- Row 20: Segunda ---
$ws.Cells.Item(20, 1).Value = 3
$ws.Cells.Item(20, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44467
$ws.Cells.Item(20, 5).Value = 5
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100107
$ws.Cells.Item(20, 8).Value = "Otros"
$ws.Cells.Item(20, 9).Value = 100107002
$ws.Cells.Item(20, 10).Value = "Chirimoya"
$ws.Cells.Item(20, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(20, 12).Value = "Segunda"
$ws.Cells.Item(20, 13).Value = 47
$ws.Cells.Item(20, 14).Value = 22000
$ws.Cells.Item(20, 15).Value = 22000
$ws.Cells.Item(20, 16).Value = 22000
$ws.Cells.Item(20, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(20, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(20, 19).Value = 2200
$ws.Cells.Item(20, 20).Value = 10
